$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 0  # H8
$ws.Cells.Item(8, 9).Value = 0  # I8
$ws.Cells.Item(8, 11).Value = 0  # K8
$ws.Cells.Item(8, 13).ClearContents()  # M8
$ws.Cells.Item(17, 8).Value = 2502350.5  # H17
$ws.Cells.Item(17, 10).Value = 2649489  # J17
$ws.Cells.Item(17, 12).Value = 7948467  # L17
$ws.Cells.Item(17, 14).Value = -7948803  # N17
$ws.Cells.Item(18, 8).Value = 429.8889  # H18
$ws.Cells.Item(18, 9).Value = 429.8889  # I18
$ws.Cells.Item(18, 11).Value = 429.8889  # K18
$ws.Cells.Item(18, 13).Value = -145.8889  # M18
$ws.Cells.Item(32, 8).Value = 8212  # H32
$ws.Cells.Item(32, 9).Value = 1750  # I32
$ws.Cells.Item(32, 11).Value = 1750  # K32
$ws.Cells.Item(32, 13).Value = -1424  # M32
$ws.Cells.Item(62, 8).Value = 5596.9473  # H62
$ws.Cells.Item(62, 9).Value = 5310.5557  # I62
$ws.Cells.Item(62, 11).Value = 5310.5557  # K62
$ws.Cells.Item(62, 13).Value = -4686.5557  # M62
$ws.Cells.Item(65, 8).Value = 5596.9473  # H65
$ws.Cells.Item(65, 9).Value = 5310.5557  # I65
$ws.Cells.Item(65, 11).Value = 26552.7785  # K65
$ws.Cells.Item(65, 13).Value = -23432.7785  # M65
$ws.Cells.Item(69, 8).Value = 0  # H69
$ws.Cells.Item(69, 10).Value = 0  # J69
$ws.Cells.Item(69, 12).Value = 0  # L69
$ws.Cells.Item(69, 14).ClearContents()  # N69
$ws.Cells.Item(72, 8).Value = 0  # H72
$ws.Cells.Item(72, 10).Value = 0  # J72
$ws.Cells.Item(72, 12).Value = 0  # L72
$ws.Cells.Item(72, 14).ClearContents()  # N72
$ws.Cells.Item(92, 8).Value = 1461.3914  # H92
$ws.Cells.Item(92, 9).Value = 1231  # I92
$ws.Cells.Item(92, 11).Value = 1231  # K92
$ws.Cells.Item(92, 13).Value = 17  # M92
$ws.Cells.Item(98, 8).Value = 62215.5  # H98
$ws.Cells.Item(98, 9).Value = 62215.5  # I98
$ws.Cells.Item(98, 11).Value = 62215.5  # K98
$ws.Cells.Item(98, 13).Value = -60717.5  # M98
$ws.Cells.Item(122, 8).Value = 62215.5  # H122
$ws.Cells.Item(122, 9).Value = 62215.5  # I122
$ws.Cells.Item(122, 11).Value = 186646.5  # K122
$ws.Cells.Item(122, 13).Value = -184196.5  # M122
$ws.Cells.Item(132, 8).Value = 1360.921  # H132
$ws.Cells.Item(132, 9).Value = 881  # I132
$ws.Cells.Item(132, 11).Value = 2643  # K132
$ws.Cells.Item(132, 13).Value = -113  # M132
$ws.Cells.Item(134, 8).Value = 193000  # H134
$ws.Cells.Item(134, 10).Value = 193000  # J134
$ws.Cells.Item(134, 12).Value = 193000  # L134
$ws.Cells.Item(134, 14).Value = -203140  # N134
$ws.Cells.Item(137, 8).Value = 80002630  # H137
$ws.Cells.Item(137, 9).Value = 45457332  # I137
$ws.Cells.Item(137, 11).Value = 136371996  # K137
$ws.Cells.Item(137, 13).Value = -136369446  # M137
$ws.Cells.Item(138, 8).Value = 3472.9333  # H138
$ws.Cells.Item(138, 9).Value = 2251  # I138
$ws.Cells.Item(138, 10).Value = 3969.3438  # J138
$ws.Cells.Item(138, 11).Value = 6753  # K138
$ws.Cells.Item(138, 12).Value = 11908.0314  # L138
$ws.Cells.Item(138, 13).Value = -1613  # M138
$ws.Cells.Item(138, 14).Value = -22188.0314  # N138

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 324.42856  # H5
$ws.Cells.Item(5, 9).Value = 154  # I5
$ws.Cells.Item(5, 11).Value = 154  # K5
$ws.Cells.Item(5, 13).Value = -42  # M5
$ws.Cells.Item(32, 8).Value = 14504325  # H32
$ws.Cells.Item(32, 9).Value = 21282104  # I32
$ws.Cells.Item(32, 11).Value = 21282104  # K32
$ws.Cells.Item(32, 13).Value = -21281817  # M32
$ws.Cells.Item(74, 8).Value = 55619836  # H74
$ws.Cells.Item(74, 9).Value = 58891396  # I74
$ws.Cells.Item(74, 11).Value = 58891396  # K74
$ws.Cells.Item(74, 13).Value = -58890522  # M74
$ws.Cells.Item(77, 8).Value = 55619836  # H77
$ws.Cells.Item(77, 9).Value = 58891396  # I77
$ws.Cells.Item(77, 11).Value = 294456980  # K77
$ws.Cells.Item(77, 13).Value = -294452612  # M77
$ws.Cells.Item(138, 8).Value = 150000  # H138
$ws.Cells.Item(138, 10).Value = 150000  # J138
$ws.Cells.Item(138, 12).Value = 150000  # L138
$ws.Cells.Item(138, 14).Value = -160280  # N138

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 324.42856  # H4
$ws.Cells.Item(4, 9).Value = 154  # I4
$ws.Cells.Item(4, 11).Value = 154  # K4
$ws.Cells.Item(4, 13).Value = -39  # M4
$ws.Cells.Item(22, 8).Value = 99.5  # H22
$ws.Cells.Item(22, 9).Value = 99  # I22
$ws.Cells.Item(22, 10).Value = 100  # J22
$ws.Cells.Item(22, 11).Value = 99  # K22
$ws.Cells.Item(22, 12).Value = 100  # L22
$ws.Cells.Item(22, 13).Value = 74  # M22
$ws.Cells.Item(22, 14).Value = -446  # N22
$ws.Cells.Item(68, 8).Value = 0  # H68
$ws.Cells.Item(68, 10).Value = 0  # J68
$ws.Cells.Item(68, 12).Value = 0  # L68
$ws.Cells.Item(68, 14).ClearContents()  # N68
$ws.Cells.Item(71, 8).Value = 0  # H71
$ws.Cells.Item(71, 10).Value = 0  # J71
$ws.Cells.Item(71, 12).Value = 0  # L71
$ws.Cells.Item(71, 14).ClearContents()  # N71
$ws.Cells.Item(94, 8).Value = 2557.1853  # H94
$ws.Cells.Item(94, 9).Value = 907.7143  # I94
$ws.Cells.Item(94, 11).Value = 907.7143  # K94
$ws.Cells.Item(94, 13).Value = -456.7143  # M94
$ws.Cells.Item(97, 8).Value = 13565.454  # H97
$ws.Cells.Item(97, 9).Value = 8922.9  # I97
$ws.Cells.Item(97, 11).Value = 8922.9  # K97
$ws.Cells.Item(97, 13).Value = -7931.9  # M97
$ws.Cells.Item(134, 8).Value = 11497.25  # H134
$ws.Cells.Item(134, 9).Value = 12000  # I134
$ws.Cells.Item(134, 10).Value = 11425.429  # J134
$ws.Cells.Item(134, 11).Value = 36000  # K134
$ws.Cells.Item(134, 12).Value = 34276.287  # L134
$ws.Cells.Item(134, 13).Value = -33465  # M134
$ws.Cells.Item(134, 14).Value = -39346.287  # N134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 12409.223  # H22
$ws.Cells.Item(22, 9).Value = 25202.5  # I22
$ws.Cells.Item(22, 10).Value = 2174.6  # J22
$ws.Cells.Item(22, 11).Value = 25202.5  # K22
$ws.Cells.Item(22, 12).Value = 2174.6  # L22
$ws.Cells.Item(22, 13).Value = -24852.5  # M22
$ws.Cells.Item(22, 14).Value = -2874.6  # N22
$ws.Cells.Item(31, 8).Value = 24394432  # H31
$ws.Cells.Item(31, 9).Value = 3295.9375  # I31
$ws.Cells.Item(31, 10).Value = 111118470  # J31
$ws.Cells.Item(31, 11).Value = 3295.9375  # K31
$ws.Cells.Item(31, 12).Value = 111118470  # L31
$ws.Cells.Item(31, 13).Value = -3000.9375  # M31
$ws.Cells.Item(31, 14).Value = -111119060  # N31
$ws.Cells.Item(34, 8).Value = 24394432  # H34
$ws.Cells.Item(34, 9).Value = 3295.9375  # I34
$ws.Cells.Item(34, 10).Value = 111118470  # J34
$ws.Cells.Item(34, 11).Value = 3295.9375  # K34
$ws.Cells.Item(34, 12).Value = 111118470  # L34
$ws.Cells.Item(34, 13).Value = -3093.9375  # M34
$ws.Cells.Item(34, 14).Value = -111118874  # N34
$ws.Cells.Item(50, 8).Value = 0  # H50
$ws.Cells.Item(50, 10).Value = 0  # J50
$ws.Cells.Item(50, 12).Value = 0  # L50
$ws.Cells.Item(50, 14).ClearContents()  # N50
$ws.Cells.Item(51, 8).Value = 47214.5  # H51
$ws.Cells.Item(51, 10).Value = 51099  # J51
$ws.Cells.Item(51, 12).Value = 51099  # L51
$ws.Cells.Item(51, 14).Value = -52571  # N51
$ws.Cells.Item(61, 8).Value = 47214.5  # H61
$ws.Cells.Item(61, 10).Value = 51099  # J61
$ws.Cells.Item(61, 12).Value = 51099  # L61
$ws.Cells.Item(61, 14).Value = -51795  # N61
$ws.Cells.Item(122, 8).Value = 2139.8333  # H122
$ws.Cells.Item(122, 9).Value = 1778.8  # I122
$ws.Cells.Item(122, 11).Value = 5336.4  # K122
$ws.Cells.Item(122, 13).Value = -2886.4  # M122

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 200  # H20
$ws.Cells.Item(20, 9).Value = 200  # I20
$ws.Cells.Item(20, 11).Value = 600  # K20
$ws.Cells.Item(20, 13).Value = -373  # M20
$ws.Cells.Item(22, 8).Value = 739.8095  # H22
$ws.Cells.Item(22, 9).Value = 274.5  # I22
$ws.Cells.Item(22, 10).Value = 2228.8  # J22
$ws.Cells.Item(22, 11).Value = 823.5  # K22
$ws.Cells.Item(22, 12).Value = 6686.400000000001  # L22
$ws.Cells.Item(22, 13).Value = -654.5  # M22
$ws.Cells.Item(22, 14).Value = -7024.400000000001  # N22
$ws.Cells.Item(27, 8).Value = 739.8095  # H27
$ws.Cells.Item(27, 9).Value = 274.5  # I27
$ws.Cells.Item(27, 10).Value = 2228.8  # J27
$ws.Cells.Item(27, 11).Value = 823.5  # K27
$ws.Cells.Item(27, 12).Value = 6686.400000000001  # L27
$ws.Cells.Item(27, 13).Value = -721.5  # M27
$ws.Cells.Item(27, 14).Value = -6890.400000000001  # N27
$ws.Cells.Item(34, 8).Value = 657.8  # H34
$ws.Cells.Item(34, 9).Value = 574.75  # I34
$ws.Cells.Item(34, 10).Value = 990  # J34
$ws.Cells.Item(34, 11).Value = 1724.25  # K34
$ws.Cells.Item(34, 12).Value = 2970  # L34
$ws.Cells.Item(34, 13).Value = -1640.25  # M34
$ws.Cells.Item(34, 14).Value = -3138  # N34
$ws.Cells.Item(39, 8).Value = 2980.9524  # H39
$ws.Cells.Item(39, 10).Value = 3792.8572  # J39
$ws.Cells.Item(39, 12).Value = 11378.5716  # L39
$ws.Cells.Item(39, 14).Value = -11966.5716  # N39
$ws.Cells.Item(131, 8).Value = 50715.523  # H131
$ws.Cells.Item(131, 10).Value = 7020.1177  # J131
$ws.Cells.Item(131, 12).Value = 21060.3531  # L131
$ws.Cells.Item(131, 14).Value = -31140.3531  # N131
$ws.Cells.Item(138, 8).Value = 1388.0834  # H138
$ws.Cells.Item(138, 9).Value = 1377.909  # I138
$ws.Cells.Item(138, 10).Value = 1500  # J138
$ws.Cells.Item(138, 11).Value = 4133.727000000001  # K138
$ws.Cells.Item(138, 12).Value = 4500  # L138
$ws.Cells.Item(138, 13).Value = 1006.272999999999  # M138
$ws.Cells.Item(138, 14).Value = -14780  # N138

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1515434  # H2
$ws.Cells.Item(2, 9).Value = 2500058.2  # I2
$ws.Cells.Item(2, 10).Value = 627.38464  # J2
$ws.Cells.Item(2, 11).Value = 2500058.2  # K2
$ws.Cells.Item(2, 12).Value = 627.38464  # L2
$ws.Cells.Item(2, 13).Value = -2499945.2  # M2
$ws.Cells.Item(2, 14).Value = -853.38464  # N2
$ws.Cells.Item(132, 8).Value = 2522.543  # H132
$ws.Cells.Item(132, 9).Value = 2266.484  # I132
$ws.Cells.Item(132, 10).Value = 4507  # J132
$ws.Cells.Item(132, 11).Value = 6799.451999999999  # K132
$ws.Cells.Item(132, 12).Value = 13521  # L132
$ws.Cells.Item(132, 13).Value = -4269.451999999999  # M132
$ws.Cells.Item(132, 14).Value = -18581  # N132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1377.9524  # H46
$ws.Cells.Item(46, 9).Value = 988.2353000000001  # I46
$ws.Cells.Item(46, 11).Value = 988.2353000000001  # K46
$ws.Cells.Item(46, 13).Value = -800.2353000000001  # M46
$ws.Cells.Item(136, 8).Value = 3834.3215  # H136
$ws.Cells.Item(136, 9).Value = 3834.3215  # I136
$ws.Cells.Item(136, 11).Value = 11502.9645  # K136
$ws.Cells.Item(136, 13).Value = -8952.9645  # M136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 0  # H49
$ws.Cells.Item(49, 9).Value = 0  # I49
$ws.Cells.Item(49, 11).Value = 0  # K49
$ws.Cells.Item(49, 13).ClearContents()  # M49
$ws.Cells.Item(81, 8).Value = 1238.0476  # H81
$ws.Cells.Item(84, 8).Value = 1238.0476  # H84
$ws.Cells.Item(132, 8).Value = 6393.0938  # H132
$ws.Cells.Item(132, 9).Value = 5961.4194  # I132
$ws.Cells.Item(132, 10).Value = 19775  # J132
$ws.Cells.Item(132, 11).Value = 17884.2582  # K132
$ws.Cells.Item(132, 12).Value = 59325  # L132
$ws.Cells.Item(132, 13).Value = -15354.2582  # M132
$ws.Cells.Item(132, 14).Value = -64385  # N132
$ws.Cells.Item(135, 8).Value = 33371666  # H135
$ws.Cells.Item(135, 9).Value = 25000  # I135
$ws.Cells.Item(135, 10).Value = 50045000  # J135
$ws.Cells.Item(135, 11).Value = 25000  # K135
$ws.Cells.Item(135, 12).Value = 50045000  # L135
$ws.Cells.Item(135, 13).Value = -19930  # M135
$ws.Cells.Item(135, 14).Value = -50055140  # N135
$ws.Cells.Item(136, 8).Value = 3650.25  # H136
$ws.Cells.Item(136, 9).Value = 2951.4  # I136
$ws.Cells.Item(136, 10).Value = 4815  # J136
$ws.Cells.Item(136, 11).Value = 8854.200000000001  # K136
$ws.Cells.Item(136, 12).Value = 14445  # L136
$ws.Cells.Item(136, 13).Value = -6304.200000000001  # M136
$ws.Cells.Item(136, 14).Value = -19545  # N136
